$d = $word.ActiveDocument

# 1. "Normal" style: change line spacing from 1.5 lines (360) to double (480),
#    i.e. match the input file's line spacing (keep the Arial font as-is).
$normal = $d.Styles.Item("Normal")
$normal.ParagraphFormat.LineSpacingRule = 2

# 2. "Heading 5" / "Heading 6" / "Title" styles: they used to pin their own
#    line spacing to double (480/auto) on top of before/after spacing; that
#    explicit override is removed so they again inherit line spacing from
#    Normal.
$heading5 = $d.Styles.Item("Heading 5")
$heading5.ParagraphFormat.LineSpacing = 0

$heading6 = $d.Styles.Item("Heading 6")
$heading6.ParagraphFormat.LineSpacing = 0

$title = $d.Styles.Item("Title")
$title.ParagraphFormat.LineSpacing = 0
